$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in column C ("Fede" review column) for the existing deliverables (rows 5-28)
$siRows = @(5,7,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,25,26,27,28)
foreach ($r in $siRows) {
    $ws.Range("C$r").Value = "si"
}

$noRows = @(6,8,24)
foreach ($r in $noRows) {
    $ws.Range("C$r").Value = "no"
}

# New deliverable rows at the bottom of the list
$ws.Range("B31").Value = "codigo fuente del front end(solo en el WBS va desagregado)"
$ws.Range("C31").Value = "si"

$ws.Range("B32").Value = "script base de datos"
$ws.Range("C32").Value = "si"

$ws.Range("B33").Value = "codigo fuente del back end"
$ws.Range("C33").Value = "si"

# Update the saved view state (scroll position + active selection)
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("C33").Select()
